# Update cryptos list values (price and 1h volume change) per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.107.66"
$ws.Range("E2").Value = "  +0.62%  "
$ws.Range("D3").Value = "2.564.76"
$ws.Range("E3").Value = "  +1.23%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "315.07"
$ws.Range("E5").Value = "  -0.92%  "
$ws.Range("D6").Value = "96.69"
$ws.Range("E6").Value = "  +0.97%  "
$ws.Range("D7").Value = "0.578"
$ws.Range("E7").Value = "  -0.41%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").Value = "0.539"
$ws.Range("E9").Value = "  +1.08%  "
$ws.Range("D10").Value = "35.43"
$ws.Range("E10").Value = "  -2.49%  "
$ws.Range("D11").Value = "0.0814"
$ws.Range("E11").Value = "  +0.32%  "
$ws.Range("D12").Value = "7.44"
$ws.Range("E12").Value = "  -2.50%  "
$ws.Range("D13").Value = "2.958.72"
$ws.Range("E13").Value = "  +1.09%  "
$ws.Range("E14").Value = "  -3.67%  "
$ws.Range("D15").Value = "2.630.24"
$ws.Range("E15").Value = "  +3.09%  "
$ws.Range("D16").Value = "15.04"
$ws.Range("E16").Value = "  -2.06%  "
$ws.Range("D17").Value = "0.842"
$ws.Range("E17").Value = "  -0.96%  "
$ws.Range("D18").Value = "43.121.41"
$ws.Range("E18").Value = "  +0.39%  "
$ws.Range("D19").Value = "6.83"
$ws.Range("E19").Value = "  +2.79%  "
$ws.Range("D20").Value = "12.58"
$ws.Range("E20").Value = "  -3.69%  "
$ws.Range("D21").Value = "0.0₃0961"
$ws.Range("E21").Value = "  -0.83%  "
$ws.Range("D22").Value = "69.26"
$ws.Range("E22").Value = "  -1.36%  "
$ws.Range("D23").Value = "253.30"
$ws.Range("E23").Value = "  +0.52%  "
$ws.Range("E24").Value = "  -0.28%  "
$ws.Range("E25").Value = "  +2.29%  "
$ws.Range("D26").Value = "26.82"
$ws.Range("E26").Value = "  -0.87%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("E28").Value = "  +0.21%  "
$ws.Range("D29").Value = "40.07"
$ws.Range("E29").Value = "  +0.44%  "
$ws.Range("D30").Value = "10.25"
$ws.Range("E30").Value = "  -0.17%  "
$ws.Range("D31").Value = "5.83"
$ws.Range("E31").Value = "  -4.12%  "
$ws.Range("D32").Value = "154.92"
$ws.Range("E32").Value = "  +0.67%  "
$ws.Range("D33").Value = "3.39"
$ws.Range("E33").Value = "  +2.11%  "
$ws.Range("E34").Value = "  +2.05%  "
$ws.Range("D35").Value = "2.71"
$ws.Range("E35").Value = "  +3.36%  "
$ws.Range("E36").Value = "  -0.27%  "
$ws.Range("D37").Value = "18.93"
$ws.Range("E37").Value = "  -0.59%  "
$ws.Range("E38").Value = "  -0.31%  "
$ws.Range("E39").Value = "  +6.70%  "
$ws.Range("D41").Value = "22.53"
$ws.Range("E41").Value = "  -5.61%  "
$ws.Range("E42").Value = "  +4.01%  "
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("E45").Value = "  -2.46%  "
$ws.Range("D46").Value = "2.002.99"
$ws.Range("E46").Value = "  -0.95%  "
$ws.Range("D47").Value = "8.87"
$ws.Range("E47").Value = "  +0.94%  "
$ws.Range("D48").Value = "2.811.32"
$ws.Range("E48").Value = "  +0.92%  "
$ws.Range("D49").Value = "82.71"
$ws.Range("E49").Value = "  -3.35%  "
$ws.Range("D50").Value = "74.84"
$ws.Range("E50").Value = "  +0.72%  "
$ws.Range("E51").Value = "  +2.24%  "
